$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the execution_datetime on the existing last row (row 14)
$ws.Cells.Item(14, 2).Value = 44382.88801283449

# New rows of data appended after row 14
$newRows = @(
    @{ A = "ytube-transcripts-text--00QUYoZHnH8.txt"; B = 44382.88804535648; C = 1; D = 16 },
    @{ A = "ytube-transcripts-text--0H5QZvOqlJM.txt"; B = 44382.88807880093; C = 3; D = 44 },
    @{ A = "ytube-transcripts-text--24R8JObNNQ4.txt"; B = 44382.88811206944; C = 2; D = 45 },
    @{ A = "ytube-transcripts-text--3HJj85K_7MQ.txt"; B = 44382.88814605093; C = 1; D = 12 }
)

$dateNumberFormat = $ws.Cells.Item(14, 2).NumberFormat

$rowIndex = 15
foreach ($row in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $row.A
    $ws.Cells.Item($rowIndex, 2).Value = $row.B
    $ws.Cells.Item($rowIndex, 2).NumberFormat = $dateNumberFormat
    $ws.Cells.Item($rowIndex, 3).Value = $row.C
    $ws.Cells.Item($rowIndex, 4).Value = $row.D
    $rowIndex++
}
